$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header cells to replace spaces/slashes with underscores
$ws.Range("E1").Value = "KW_BW"
$ws.Range("F1").Value = "Heart_Weight"
$ws.Range("G1").Value = "HW_BW"

# Update the active selection to F2 (was A2)
$ws.Range("F2").Select()
